$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 101, shifting existing rows 101:205 down to 102:206.
$ws.Rows(101).Insert()

# Populate the newly inserted row 101 with the new record.
$ws.Range("A101").Value = 11
$ws.Range("B101").Value = "Vega Monumental Concepción"
$ws.Range("C101").Value = "Bíobío"
$ws.Range("D101").Value = 44810
$ws.Range("E101").Value = 8
$ws.Range("F101").Value = 100112003
$ws.Range("G101").Value = "Ajo"
$ws.Range("H101").Value = "Chino"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 400
$ws.Range("K101").Value = 17000
$ws.Range("L101").Value = 18000
$ws.Range("M101").Value = 17500
$ws.Range("N101").Value = "$/caja 10 kilos"
$ws.Range("O101").Value = "China"
$ws.Range("P101").Value = 1750
$ws.Range("Q101").Value = 10
$ws.Range("R101").Value = "Hortaliza"
